$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at the left (A and B), shifting existing data right.
$ws.Range("A1:B6").Insert(-4161)

# The old "codigo_venta" column (originally P, now shifted to R) moves to
# the new column B, right after the new "nro_venta" column.
$ws.Range("R1:R6").Cut($ws.Range("B1:B6"))

# Populate the brand-new "nro_venta" column (A).
$ws.Range("A1").Value = "nro_venta"
$ws.Range("A2").Value = 7
$ws.Range("A3").Value = 8
$ws.Range("A4").Value = 9
$ws.Range("A5").Value = 10
$ws.Range("A6").Value = 11

# Match the header formatting (bold/centered/bordered) used by the rest of row 1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear out the now-empty column R left behind by the cut.
$ws.Range("R1:R6").Clear()
